# Generate Report for Handoff
# Adds two new handed-off files (ad53d21f..., f5d6b03e...) to every sheet of the
# localization status report: the "Overview" rollup sheet and the per-locale
# "zh-cn" / "de-de" detail sheets. Mirrors the existing "Ready for handoff"
# rows (e.g. row 3, file 462da0bf...) already present in each sheet.

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Overview sheet — one row per new file, columns A..G
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$overviewRows = @(
    @{ Row = 4; Name = "ad53d21f-7296-4bd1-8b5a-4571ffd30a4f.md"; Date = "2016-08-22 18:43:45" },
    @{ Row = 5; Name = "f5d6b03e-4942-4da3-b963-c9b6cd9c01c7.md"; Date = "2016-08-22 18:43:45" }
)

foreach ($r in $overviewRows) {
    $row = $r.Row
    $fileName = $r.Name
    $pathName = "e2e\" + $fileName

    $ws.Range("A" + $row).Value = $fileName
    $ws.Range("B" + $row).Value = $pathName
    $ws.Range("C" + $row).Value = ".md"
    $ws.Range("D" + $row).Value = ""
    $ws.Range("E" + $row).Value = "Ready for handoff"
    $ws.Range("F" + $row).Value = "Ready for handoff"
    $ws.Range("G" + $row).Value = $r.Date
    $ws.Range("G" + $row).NumberFormat = $dateFmt

    $ws.Hyperlinks.Add($ws.Range("B" + $row), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/" + $fileName, "", "", $pathName) | Out-Null
    $ws.Range("B" + $row).Style = "HyperLink"
}

$loOverview = $ws.ListObjects.Item(1)
$loOverview.Resize($ws.Range("A1:G5"))

# ---------------------------------------------------------------------------
# Per-locale detail sheets — one row per new file, columns A..P
# ---------------------------------------------------------------------------
$localeSheets = @(
    @{ Sheet = "zh-cn"; Locale = "zh-cn" },
    @{ Sheet = "de-de"; Locale = "de-de" }
)

$detailRows = @(
    @{ Row = 4; Name = "ad53d21f-7296-4bd1-8b5a-4571ffd30a4f.md"; Hash = "a6ca21533bf6e5d8c0e6d758a7357f3fbb8304c1" },
    @{ Row = 5; Name = "f5d6b03e-4942-4da3-b963-c9b6cd9c01c7.md"; Hash = "21bd2561cda9d3f1b346f8f49e80bf7b5da938c7" }
)

# zh-cn handoff xliffs were generated a few minutes before the de-de ones
$handoffDates = @{ "zh-cn" = "2016-08-22 18:43:40"; "de-de" = "2016-08-22 18:43:45" }

foreach ($loc in $localeSheets) {
    $ws = $wb.Worksheets.Item($loc.Sheet)
    $locale = $loc.Locale
    $handoffDate = $handoffDates[$locale]

    foreach ($r in $detailRows) {
        $row = $r.Row
        $fileName = $r.Name
        $baseName = $fileName.Substring(0, $fileName.Length - 3)  # strip trailing ".md"
        $xliffName = $baseName + "." + $r.Hash + "." + $locale + ".xlf"

        $ws.Range("A" + $row).Value = $fileName
        $ws.Range("B" + $row).Value = ".md"
        $ws.Range("C" + $row).Value = "Ready for handoff"
        $ws.Range("D" + $row).Value = "e2e"
        $ws.Range("E" + $row).Value = "ht"
        $ws.Range("F" + $row).Value = "False"
        $ws.Range("G" + $row).Value = $xliffName
        $ws.Range("H" + $row).Value = $handoffDate
        $ws.Range("H" + $row).NumberFormat = $dateFmt
        $ws.Range("I" + $row).Value = ""
        $ws.Range("J" + $row).Value = ""
        $ws.Range("K" + $row).Value = "0001-01-01 00:00:00"
        $ws.Range("K" + $row).NumberFormat = $dateFmt
        $ws.Range("L" + $row).Value = ""
        $ws.Range("M" + $row).Value = "True"
        $ws.Range("N" + $row).Value = ""
        $ws.Range("O" + $row).Value = "False"
        $ws.Range("P" + $row).Value = ""

        $ws.Hyperlinks.Add($ws.Range("A" + $row), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/" + $fileName, "", "", $fileName) | Out-Null
        $ws.Range("A" + $row).Style = "HyperLink"
    }

    $lo = $ws.ListObjects.Item(1)
    $lo.Resize($ws.Range("A1:P5"))
}
